$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Tighten up the "Attribute" strings: drop the space after the semicolon
# that separates attribute definitions (cosmetic text cleanup).
$ws.Range("D3").Value = "Ticketart:Dropdown(Bus,Zug,U-Bahn);Häufigkeit:Dropdown(Täglich,Wöchentlich,Selten)"
$ws.Range("D6").Value = "Land:Text;Häufigkeit:Dropdown(Selten,Oft,Regelmäßig)"
$ws.Range("D7").Value = "Anteil:Number;Technik:Dropdown(Laptop,PC,Tablet)"

# These rows no longer need their tall custom row height now that the
# text is shorter - let Excel recompute the natural (auto) row height.
$ws.Rows.Item(2).AutoFit()
$ws.Rows.Item(4).AutoFit()
$ws.Rows.Item(7).AutoFit()

# Move the active selection from E9 to D9.
$ws.Range("D9").Select()
